# Atualização de bases das ligas, do dia: 21-04-2024 às 13:33
#
# This script applies the changes described by the diff:
#  1) Four pairs of rows had their HomeTeam/AwayTeam/odds data (columns B:AC)
#     swapped between the two rows (the "id" in column A stays tied to the
#     row position, not the match).
#  2) Two "upcoming match" rows near the bottom of the sheet were removed
#     (matches that no longer needed to be tracked), which shifted the
#     remaining upcoming-match rows up, and their odds were refreshed with
#     newer snapshot values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-Rows($r1, $r2) {
    # Columns B..AC are column indices 2..29; column A (id) is left untouched.
    for ($col = 2; $col -le 29; $col++) {
        $c1 = $ws.Cells.Item($r1, $col)
        $c2 = $ws.Cells.Item($r2, $col)
        $v1 = $c1.Value2
        $v2 = $c2.Value2
        $c1.Value = $v2
        $c2.Value = $v1
    }
}

# --- 1) Swap the four row pairs ---
Swap-Rows 95 96
Swap-Rows 110 111
Swap-Rows 129 131
Swap-Rows 237 238

# --- 2) Remove the two stale upcoming-match rows (old rows 242 & 243) ---
$ws.Range("A242:A243").EntireRow.Delete()

# After the delete, old rows 244-247 have shifted up to become rows 242-245.
# Column A (the sequential id) must stay tied to row position, so restore it.
$ws.Cells.Item(242, 1).Value = 240
$ws.Cells.Item(243, 1).Value = 241
$ws.Cells.Item(244, 1).Value = 242
$ws.Cells.Item(245, 1).Value = 243

# Refresh the odds snapshots for the newly-shifted rows to the latest values.
# Row 242 (AD Guanacasteca vs Herediano)
$ws.Range("R242").Value = 1.9
$ws.Range("S242").Value = 1.9

# Row 243 (Santos de Gupiles vs Puntarenas)
$ws.Range("N243").Value = 2.15
$ws.Range("O243").Value = 3.25
$ws.Range("P243").Value = 3.1
$ws.Range("Q243").Value = -0.25
$ws.Range("R243").Value = 1.9
$ws.Range("S243").Value = 1.9
$ws.Range("T243").Value = 2.25
$ws.Range("U243").Value = 2.05
$ws.Range("V243").Value = 1.75

# Row 244 (Municipal Liberia vs Cartagines)
$ws.Range("P244").Value = 3
$ws.Range("U244").Value = 1.95
$ws.Range("V244").Value = 1.85

# Row 245 (AD San Carlos vs Sporting San Jose)
$ws.Range("N245").Value = 1.4
$ws.Range("O245").Value = 4
$ws.Range("P245").Value = 6
$ws.Range("Q245").Value = -1.25
$ws.Range("R245").Value = 1.975
$ws.Range("S245").Value = 1.825
$ws.Range("U245").Value = 1.775
$ws.Range("V245").Value = 2.025
